# Refresh cryptos list (price + volume(1h) %) to latest scraped values.
# Some "Price" cells are plain decimals (e.g. "579.61") which Excel would
# otherwise auto-convert to a Number; force those to stay Text so the
# cell type matches the rest of the (string-typed) price column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.253.00'
$ws.Range('E2').Value = '  -4.49%  '
$ws.Range('D3').Value = '2.988.64'
$ws.Range('E3').Value = '  -5.91%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.61'
$ws.Range('E5').Value = '  -1.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '125.69'
$ws.Range('E6').Value = '  -6.80%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '2.987.28'
$ws.Range('E8').Value = '  -5.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.498'
$ws.Range('E9').Value = '  -3.15%  '
$ws.Range('E10').Value = '  -6.00%  '
$ws.Range('E11').Value = '  -2.08%  '
$ws.Range('E12').Value = '  -2.97%  '
$ws.Range('E13').Value = '  -5.72%  '
$ws.Range('E14').Value = '  -5.33%  '
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').Value = '3.474.73'
$ws.Range('E16').Value = '  -6.10%  '
$ws.Range('D17').Value = '2.977.80'
$ws.Range('E17').Value = '  -6.20%  '
$ws.Range('D18').Value = '60.190.31'
$ws.Range('E18').Value = '  -4.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.22'
$ws.Range('E19').Value = '  -4.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '430.77'
$ws.Range('E20').Value = '  -6.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.06'
$ws.Range('E21').Value = '  -6.37%  '
$ws.Range('E22').Value = '  -4.89%  '
$ws.Range('E23').Value = '  -7.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.75'
$ws.Range('E24').Value = '  -3.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '79.02'
$ws.Range('E25').Value = '  -4.07%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.56'
$ws.Range('E28').Value = '  -3.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.27'
$ws.Range('E29').Value = '  -5.03%  '
$ws.Range('E30').Value = '  -6.92%  '
$ws.Range('E31').Value = '  -8.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.34'
$ws.Range('E32').Value = '  -6.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0935'
$ws.Range('E33').Value = '  -7.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.16'
$ws.Range('E34').Value = '  -8.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.955'
$ws.Range('E35').Value = '  -7.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.58'
$ws.Range('E36').Value = '  -3.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '49.53'
$ws.Range('E37').Value = '  -3.12%  '
$ws.Range('E38').Value = '  -6.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.06'
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('E40').Value = '  -6.78%  '
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '383.04'
$ws.Range('E42').Value = '  -4.51%  '
$ws.Range('E43').Value = '  -7.00%  '
$ws.Range('D44').Value = '2.635.69'
$ws.Range('E44').Value = '  -6.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.236'
$ws.Range('E46').Value = '  -6.07%  '
$ws.Range('E47').Value = '  -5.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '118.38'
$ws.Range('E48').Value = '  -5.32%  '
$ws.Range('E49').Value = '  -3.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.54'
$ws.Range('E50').Value = '  -6.18%  '
$ws.Range('E51').Value = '  -5.38%  '
